$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 290, shifting existing rows 290-300 down to 291-301.
$ws.Rows.Item(290).Insert()

# Populate the newly inserted row 290 with its values.
$ws.Cells.Item(290, 1).Value = 4
$ws.Cells.Item(290, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(290, 3).Value = "Los Lagos"
$ws.Cells.Item(290, 4).Value = 44939
$ws.Cells.Item(290, 4).NumberFormat = $ws.Cells.Item(291, 4).NumberFormat
$ws.Cells.Item(290, 5).Value = 10
$ws.Cells.Item(290, 6).Value = "Fruta"
$ws.Cells.Item(290, 7).Value = 100108
$ws.Cells.Item(290, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(290, 9).Value = 100108002
$ws.Cells.Item(290, 10).Value = "Mango"
$ws.Cells.Item(290, 11).Value = "Sin especificar"
$ws.Cells.Item(290, 12).Value = "Primera"
$ws.Cells.Item(290, 13).Value = 200
$ws.Cells.Item(290, 14).Value = 7500
$ws.Cells.Item(290, 15).Value = 8000
$ws.Cells.Item(290, 16).Value = 7750
$ws.Cells.Item(290, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(290, 18).Value = "Brasil"
$ws.Cells.Item(290, 19).Value = 1938
$ws.Cells.Item(290, 20).Value = 4
